$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Home win": update rows 2-3 with latest predictions, drop rows 4-7
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Home win")

$ws.Cells.Item(2,1).Value = "29-12-2024 15:00"
$ws.Cells.Item(2,2).Value = "ENGLAND"
$ws.Cells.Item(2,3).Value = "LEAGUE ONE"
$ws.Cells.Item(2,4).Value = "Reading - Mansfield Town"
$ws.Cells.Item(2,5).Value = 70
$ws.Cells.Item(2,6).Value = 2.1

$ws.Cells.Item(3,1).Value = "29-12-2024 12:30"
$ws.Cells.Item(3,2).Value = "ENGLAND"
$ws.Cells.Item(3,3).Value = "LEAGUE ONE"
$ws.Cells.Item(3,4).Value = "Rotherham - Stockport County"
$ws.Cells.Item(3,5).Value = 80
$ws.Cells.Item(3,6).Value = 2.8

$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet "Draw": update rows 2-4, add new row 5
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Draw")

$ws.Cells.Item(2,1).Value = "28-12-2024 11:30"
$ws.Cells.Item(2,2).Value = "IRAN"
$ws.Cells.Item(2,3).Value = "AZADEGAN LEAGUE"
$ws.Cells.Item(2,4).Value = "Naft Gachsaran - Mes Soongoun"
$ws.Cells.Item(2,5).Value = 60
$ws.Cells.Item(2,6).Value = 2.85

$ws.Cells.Item(3,1).Value = "28-12-2024 14:30"
$ws.Cells.Item(3,2).Value = "WORLD"
$ws.Cells.Item(3,3).Value = "GULF CUP OF NATIONS"
$ws.Cells.Item(3,4).Value = "Bahrain - Yemen"
$ws.Cells.Item(3,5).Value = 60
$ws.Cells.Item(3,6).Value = 4.75

$ws.Cells.Item(4,1).Value = "29-12-2024 14:00"
$ws.Cells.Item(4,2).Value = "ITALY"
$ws.Cells.Item(4,3).Value = "SERIE B"
$ws.Cells.Item(4,4).Value = "Bari - Spezia"
$ws.Cells.Item(4,5).Value = 60
$ws.Cells.Item(4,6).Value = 2.85

$ws.Cells.Item(5,1).Value = "29-12-2024 12:30"
$ws.Cells.Item(5,2).Value = "CAMEROON"
$ws.Cells.Item(5,3).Value = "ELITE ONE"
$ws.Cells.Item(5,4).Value = "Bamboutos - Fauve Azur Elite"
$ws.Cells.Item(5,5).Value = 60
$ws.Cells.Item(5,6).Value = 2.9

# ---------------------------------------------------------------------
# Sheet "Btts": update rows 2-4, add new rows 5-8
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Btts")

$ws.Cells.Item(2,1).Value = "28-12-2024 15:00"
$ws.Cells.Item(2,2).Value = "SCOTLAND"
$ws.Cells.Item(2,3).Value = "CHAMPIONSHIP"
$ws.Cells.Item(2,4).Value = "Partick - Queen's Park"
$ws.Cells.Item(2,5).Value = 76.7
$ws.Cells.Item(2,6).Value = 1.8

$ws.Cells.Item(3,1).Value = "28-12-2024 18:00"
$ws.Cells.Item(3,2).Value = "ISRAEL"
$ws.Cells.Item(3,3).Value = "STATE CUP"
$ws.Cells.Item(3,4).Value = "Maccabi Tel Aviv - Hapoel Katamon"
$ws.Cells.Item(3,5).Value = 80
$ws.Cells.Item(3,6).Value = 1.77

$ws.Cells.Item(4,1).Value = "28-12-2024 16:00"
$ws.Cells.Item(4,2).Value = "WORLD"
$ws.Cells.Item(4,3).Value = "AFRICAN NATIONS CHAMPIONSHIP - QUALIFICATION"
$ws.Cells.Item(4,4).Value = "Cameroon - Central African Republic"
$ws.Cells.Item(4,5).Value = 76
$ws.Cells.Item(4,6).Value = 2.3

$ws.Cells.Item(5,1).Value = "28-12-2024 16:00"
$ws.Cells.Item(5,2).Value = "WORLD"
$ws.Cells.Item(5,3).Value = "AFRICAN NATIONS CHAMPIONSHIP - QUALIFICATION"
$ws.Cells.Item(5,4).Value = "Nigeria - Ghana"
$ws.Cells.Item(5,5).Value = 76.7
$ws.Cells.Item(5,6).Value = 2.15

$ws.Cells.Item(6,1).Value = "29-12-2024 15:00"
$ws.Cells.Item(6,2).Value = "ENGLAND"
$ws.Cells.Item(6,3).Value = "LEAGUE ONE"
$ws.Cells.Item(6,4).Value = "Exeter City - Crawley Town"
$ws.Cells.Item(6,5).Value = 80
$ws.Cells.Item(6,6).Value = 1.8

$ws.Cells.Item(7,1).Value = "29-12-2024 15:00"
$ws.Cells.Item(7,2).Value = "SCOTLAND"
$ws.Cells.Item(7,3).Value = "PREMIERSHIP"
$ws.Cells.Item(7,4).Value = "Motherwell - Rangers"
$ws.Cells.Item(7,5).Value = 76.7
$ws.Cells.Item(7,6).Value = 1.95

$ws.Cells.Item(8,1).Value = "29-12-2024 15:30"
$ws.Cells.Item(8,2).Value = "PORTUGAL"
$ws.Cells.Item(8,3).Value = "PRIMEIRA LIGA"
$ws.Cells.Item(8,4).Value = "Rio Ave - Nacional"
$ws.Cells.Item(8,5).Value = 86.7
$ws.Cells.Item(8,6).Value = 1.8

# ---------------------------------------------------------------------
# Sheet "Over_Under": overwrite rows 2-6 with latest predictions, drop row 7
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Over_Under")

$ws.Cells.Item(2,1).Value = "28-12-2024 15:00"
$ws.Cells.Item(2,2).Value = "SCOTLAND"
$ws.Cells.Item(2,3).Value = "LEAGUE ONE"
$ws.Cells.Item(2,4).Value = "Dumbarton - Kelty Hearts"
$ws.Cells.Item(2,5).Value = 85
$ws.Cells.Item(2,6).Value = 1.67
$ws.Cells.Item(2,7).Value = 60
$ws.Cells.Item(2,8).Value = 2.65

$ws.Cells.Item(3,1).Value = "28-12-2024 15:00"
$ws.Cells.Item(3,2).Value = "SCOTLAND"
$ws.Cells.Item(3,3).Value = "LEAGUE TWO"
$ws.Cells.Item(3,4).Value = "Spartans - Edinburgh City"
$ws.Cells.Item(3,5).Value = 80
$ws.Cells.Item(3,6).Value = 1.91
$ws.Cells.Item(3,7).Value = 46.7
$ws.Cells.Item(3,8).Value = 3.1

$ws.Cells.Item(4,1).Value = "28-12-2024 06:00"
$ws.Cells.Item(4,2).Value = "AUSTRALIA"
$ws.Cells.Item(4,3).Value = "A-LEAGUE"
$ws.Cells.Item(4,4).Value = "Central Coast Mariners - Auckland"
$ws.Cells.Item(4,5).Value = 86.7
$ws.Cells.Item(4,6).Value = 2
$ws.Cells.Item(4,7).Value = 53.3
$ws.Cells.Item(4,8).Value = 3.4

$ws.Cells.Item(5,1).Value = "29-12-2024 15:00"
$ws.Cells.Item(5,2).Value = "ENGLAND"
$ws.Cells.Item(5,3).Value = "CHAMPIONSHIP"
$ws.Cells.Item(5,4).Value = "Oxford United - Plymouth"
$ws.Cells.Item(5,5).Value = 60
$ws.Cells.Item(5,6).Value = 1.73
$ws.Cells.Item(5,7).Value = 60
$ws.Cells.Item(5,8).Value = 2.75

$ws.Cells.Item(6,1).Value = "29-12-2024 15:00"
$ws.Cells.Item(6,2).Value = "ENGLAND"
$ws.Cells.Item(6,3).Value = "LEAGUE ONE"
$ws.Cells.Item(6,4).Value = "Exeter City - Crawley Town"
$ws.Cells.Item(6,5).Value = 80
$ws.Cells.Item(6,6).Value = 1.9
$ws.Cells.Item(6,7).Value = 30
$ws.Cells.Item(6,8).Value = 3.25

$ws.Rows.Item(7).Delete()

# "Away Win" sheet is unchanged by this update.
